# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Input Parameters")
$ws2 = $wb.Worksheets.Item("Output Results")

# --- Sheet 1 ("Input Parameters"): normalise Start/End Date columns (D,E) for rows 2-7 ---
for ($r = 2; $r -le 7; $r++) {
    $ws1.Cells.Item($r, 4).Value = 41275
    $ws1.Cells.Item($r, 5).Value = 42735
}

# --- Sheet 2 ("Output Results"): rewrite rows 2-9 and append new rows 10-25 ---
$data = @(
    @("A1", 0, "Maize", 41527, 251,  14.17398169766591, 151.2598850610039),
    @("A1", 1, "Maize", 41892, 616,  14.30457553090168, 130.2603605992372),
    @("A1", 2, "Maize", 42257, 981,  14.10789703038272, 151.9922915062926),
    @("A1", 3, "Maize", 42623, 1347, 14.18220979331418, 123.6655511553475),
    @("A2", 0, "Maize", 41527, 251,  14.17398169766591, 170.8047151792224),
    @("A2", 1, "Maize", 41892, 616,  14.30457553090168, 150.0545770031289),
    @("A2", 2, "Maize", 42257, 981,  14.10789703038272, 169.1330952719165),
    @("A2", 3, "Maize", 42623, 1347, 14.18220979331418, 145.5912317627781),
    @("B1", 0, "Maize", 41543, 267,  14.15209118744616, 110.7897022513593),
    @("B1", 1, "Maize", 41908, 632,  14.11859737863568, 130.7705436744614),
    @("B1", 2, "Maize", 42273, 997,  13.76502684284506, 156.4679018704647),
    @("B1", 3, "Maize", 42639, 1363, 14.23680958115403, 136.4295070208339),
    @("B2", 0, "Maize", 41543, 267,  14.15209118744616, 128.3873177627219),
    @("B2", 1, "Maize", 41908, 632,  14.11859737863568, 148.8242877716196),
    @("B2", 2, "Maize", 42273, 997,  13.76502684284506, 174.1759787517644),
    @("B2", 3, "Maize", 42639, 1363, 14.23680958115403, 154.0218254956336),
    @("C1", 0, "Maize", 41527, 251,  12.80812688082464, 0),
    @("C1", 1, "Maize", 41892, 616,  14.20693048623959, 0),
    @("C1", 2, "Maize", 42257, 981,  13.83671734355427, 0),
    @("C1", 3, "Maize", 42623, 1347, 14.09061648046096, 0),
    @("C2", 0, "Maize", 41527, 251,  12.24033654832572, 0),
    @("C2", 1, "Maize", 41892, 616,  13.7063128475466,  0),
    @("C2", 2, "Maize", 42257, 981,  13.61443907375126, 0),
    @("C2", 3, "Maize", 42623, 1347, 13.8058694166582,  0)
)

# Number format used by the existing "Harvest Date" column (column D)
$dateFormat = $ws2.Cells.Item(2, 4).NumberFormat

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws2.Cells.Item($row, 1).Value = $vals[0]
    $ws2.Cells.Item($row, 2).Value = $vals[1]
    $ws2.Cells.Item($row, 3).Value = $vals[2]
    $ws2.Cells.Item($row, 4).Value = $vals[3]
    $ws2.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws2.Cells.Item($row, 5).Value = $vals[4]
    $ws2.Cells.Item($row, 6).Value = $vals[5]
    $ws2.Cells.Item($row, 7).Value = $vals[6]
}
